$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 5 data
$ws.Range("A5").Value = "D001"

# Set the number format before assigning the value so Excel applies the
# custom date format (dd/MM/yy) directly instead of creating an extra
# auto-generated numFmt entry.
$ws.Range("B5").NumberFormat = "dd/MM/yy"
$ws.Range("B5").Value = "10/28/2024"

$ws.Range("C5").Value = "10:00 AM"
$ws.Range("D5").Value = "12:00 PM"
$ws.Range("E5").Value = "BUSY"

# Update selection to match target state
$ws.Range("A5:E6").Select()
